# Adds a new "2022-Q4" quarter: a summary row on the "总计" sheet plus a
# brand-new "2022-Q4" detail sheet (inserted right after "总计"), which
# pushes every existing quarter sheet one tab later.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and renumber
#    the leading index column (A) for the rows that shift down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Reuse the existing index-column formatting (border/bold/centering) for
# the newly inserted A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.3

# Renumber column A sequentially (0..6) for the rows pushed down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q4" detail sheet, inserted right after "总计" (so it
#    becomes the 2nd sheet, ahead of the former "2022-Q3" sheet).
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# The "2022-Q3" tab (shifted to position 3 by the insert above) carries
# the header/index-column formatting we want to replicate.
$templateSheet = $wb.Worksheets.Item(3)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$templateSheet.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# Row 2 - 001092 广发纳斯达克生物科技指数人民币A
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "001092"
$q4.Range("C2").Value = "广发纳斯达克生物科技指数人民币A"
Set-TextValue $q4.Range("D2") "1.80"
Set-TextValue $q4.Range("E2") "90.20"
Set-TextValue $q4.Range("F2") "6.75"
$q4.Range("G2").Value = 0.1215
$q4.Range("H2").Value = 2

# Row 3 - 001093 广发纳斯达克生物科技指数美元A
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "001093"
$q4.Range("C3").Value = "广发纳斯达克生物科技指数美元A"
Set-TextValue $q4.Range("D3") "1.80"
Set-TextValue $q4.Range("E3") "90.20"
Set-TextValue $q4.Range("F3") "6.75"
$q4.Range("G3").Value = 0.1215
$q4.Range("H3").Value = 2

# Row 4 - 513290 汇添富纳斯达克生物科技ETF（QDII）
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "513290"
$q4.Range("C4").Value = "汇添富纳斯达克生物科技ETF（QDII）"
Set-TextValue $q4.Range("D4") "0.75"
Set-TextValue $q4.Range("E4") "99.38"
Set-TextValue $q4.Range("F4") "7.51"
$q4.Range("G4").Value = 0.0563
$q4.Range("H4").Value = 2

# Row 5 - 016470 广发纳斯达克生物科技指数人民币C
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "016470"
$q4.Range("C5").Value = "广发纳斯达克生物科技指数人民币C"
Set-TextValue $q4.Range("D5") "0.00"
Set-TextValue $q4.Range("E5") "90.20"
Set-TextValue $q4.Range("F5") "6.75"
$q4.Range("G5").Value = 0
$q4.Range("H5").Value = 2

# Row 6 - 016471 广发纳斯达克生物科技指数美元C
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "016471"
$q4.Range("C6").Value = "广发纳斯达克生物科技指数美元C"
Set-TextValue $q4.Range("D6") "0.00"
Set-TextValue $q4.Range("E6") "90.20"
Set-TextValue $q4.Range("F6") "6.75"
$q4.Range("G6").Value = 0
$q4.Range("H6").Value = 2

$templateSheet.Range("A2:A6").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122)

# Restore the index-column values PasteSpecial(Formats) leaves untouched
# (PasteSpecial with Formats only copies formatting, not values, so the
# numbers set above remain - nothing further to do here).

$q4.Range("A1").Select()
